# Update unit testing for database create method
#
# "Database Create" sheet previously had a single test-case row
# ("Creates Database" / "No Exceptions (void)"). It is replaced with two
# rows that describe the two real partitions for the create-database
# method: the database not existing yet (succeeds -> TRUE) and the
# database already existing (fails -> FALSE).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Database Create")

# Row 2: valid partition - database doesn't exist yet, create succeeds.
$ws.Range("C2").Value = "Database doesn't exist"
$ws.Range("D2").Value = "Database doesn't exist"
$ws.Range("E2").Value = $true

# Row 3 (new): invalid/duplicate partition - database already exists.
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Valid"
$ws.Range("C3").Value = "Database already exists"
$ws.Range("D3").Value = "Database already exists"
$ws.Range("E3").Value = $false

# Widen the Partition/Test Inputs columns to fit the longer text.
$ws.Columns.Item(3).ColumnWidth = 19.166666666666668
$ws.Columns.Item(4).ColumnWidth = 19.166666666666668

# This sheet becomes the active tab / selected cell instead of
# "Database ProcessQuery".
$ws.Activate() | Out-Null
$ws.Range("C6").Select() | Out-Null
